# REVER_DailyTracker_MONISHA.xlsx - "Add files via upload" edit
# Target sheet is OCT-2020 (the active/tab-selected sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 30 (29 Oct 2020): fill in the previously-blank task row ----
$ws.Range("C30").Value = "QMVAR"
$ws.Range("D30").Value = "QMVAR -issues fixing"
$ws.Range("E30").Value = 0.9
$ws.Range("F30").Value = "Completed"

# ---- Row 31 (30 Oct 2020): fill in the previously-blank task row ----
$ws.Range("C31").Value = "QMVAR"
$ws.Range("D31").Value = "QMVAR -Newly highlidhted  issues fixing"
$ws.Range("E31").Value = 0.8
$ws.Range("F31").Value = "Completed"

# ---- Row 32 (new row, 31 Oct 2020 - Week off) ----
# Copy the cell formatting from existing cells that already carry the
# exact styles this new row needs, then overwrite the values/text.
$ws.Range("A28").Copy()
$ws.Range("A32").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B28").Copy()
$ws.Range("B32").PasteSpecial(-4122)

$ws.Range("C28").Copy()
$ws.Range("C32").PasteSpecial(-4122)

$ws.Range("D3").Copy()
$ws.Range("D32").PasteSpecial(-4122)

$ws.Range("E28").Copy()
$ws.Range("E32").PasteSpecial(-4122)

$ws.Range("G28").Copy()
$ws.Range("F32").PasteSpecial(-4122)

$ws.Range("G28").Copy()
$ws.Range("G32").PasteSpecial(-4122)

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 44135
$ws.Range("D32").Value = "Week off"

# ---- Move the saved selection to E35, matching the new cursor position ----
$ws.Range("E35").Select()
